# Apply cryptos.xlsx price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts them away from the
# original inline-string representation (e.g. "1.000" -> 1).
$textRefs = @(
    'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13',
    'D14', 'D15', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D26',
    'D27', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38',
    'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49',
    'D50', 'D51'
)
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = '@'
}

# --- cell value updates ---
$ws.Range('D2').Value = '28.248.40'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '1.793.65'
$ws.Range('E3').Value = '  +1.82%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '326.72'
$ws.Range('E5').Value = '  -2.61%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.4462'
$ws.Range('E7').Value = '  +13.63%  '
$ws.Range('D8').Value = '0.3742'
$ws.Range('E8').Value = '  +10.22%  '
$ws.Range('D9').Value = '44.65'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').Value = '1.144'
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('D11').Value = '0.07514'
$ws.Range('E11').Value = '  +3.90%  '
$ws.Range('D12').Value = '22.61'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').Value = '1.002'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').Value = '6.279'
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('D15').Value = '7.531'
$ws.Range('E15').Value = '  +5.97%  '
$ws.Range('D16').Value = '1.783.42'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').Value = '0.00001088'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').Value = '0.06726'
$ws.Range('E18').Value = '  +1.61%  '
$ws.Range('D19').Value = '80.97'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '17.53'
$ws.Range('E21').Value = '  +3.24%  '
$ws.Range('D22').Value = '6.326'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('D23').Value = '28.234.37'
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').Value = '11.74'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('D26').Value = '20.43'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').Value = '151.90'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('D29').Value = '1.988.01'
$ws.Range('E29').Value = '  +3.24%  '
$ws.Range('D30').Value = '132.66'
$ws.Range('E30').Value = '  +2.53%  '
$ws.Range('D31').Value = '1.223'
$ws.Range('E31').Value = '  -4.22%  '
$ws.Range('D32').Value = '4.023'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').Value = '5.801'
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').Value = '0.09397'
$ws.Range('E34').Value = '  +7.60%  '
$ws.Range('D35').Value = '0.2329'
$ws.Range('E35').Value = '  +10.13%  '
$ws.Range('D36').Value = '12.09'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').Value = '0.06330'
$ws.Range('E37').Value = '  +2.24%  '
$ws.Range('D38').Value = '0.02326'
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('D39').Value = '5.166'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Value = '0.6541'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').Value = '8.294'
$ws.Range('E41').Value = '  +4.97%  '
$ws.Range('D42').Value = '1.469'
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '14.08'
$ws.Range('E45').Value = '  +2.31%  '
$ws.Range('D46').Value = '0.6084'
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('D47').Value = '3.782'
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('D48').Value = '129.92'
$ws.Range('E48').Value = '  +2.42%  '
$ws.Range('D49').Value = '2.024'
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D50').Value = '0.07122'
$ws.Range('E50').Value = '  +1.66%  '
$ws.Range('D51').Value = '1.158'
$ws.Range('E51').Value = '  -0.02%  '

# Restore the default style on the forced-text cells so only their
# content differs from the original workbook (no stray style refs).
foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = 'Normal'
}
